$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column R ("backup") ---
$ws.Range("R1").Value = "backup"
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)   # xlPasteFormats - copy header style/border

# Fill R2:R63 with 0 (existing data rows)
$ws.Range("R2:R63").Value = 0

# --- Corrections to existing rows ---
$ws.Range("Q38").Value = 0
$ws.Range("O63").Value = 2

# --- New rows 64-69 ---
$newRows = @(
    @{Row=64; A=45474; B=1645.556954586147;  C=1923.506435127782;  D=1631.594653153294;  E=1837.638305664062; G=33195846; H=2024; I=7;  J=1; K=0; L=0; M=0; N=27; O=1; P=0; Q=0},
    @{Row=65; A=45505; B=1868.455058400941;  C=1915.827150887754;  D=1680.462649269061;  E=1810.26220703125;  G=15790629; H=2024; I=8;  J=1; K=0; L=0; M=0; N=31; O=0; P=0; Q=0},
    @{Row=66; A=45536; B=1815.150024414062;  C=1859.949951171875;  D=1625;                E=1628;               G=17167543; H=2024; I=9;  J=1; K=0; L=0; M=0; N=35; O=0; P=0; Q=0},
    @{Row=67; A=45566; B=1641.5;              C=1822.900024414062;  D=1309.050048828125;  E=1393.25;            G=35510429; H=2024; I=10; J=1; K=0; L=0; M=0; N=40; O=0; P=0; Q=1},
    @{Row=68; A=45597; B=1425;                C=1520;                D=1283.25;            E=1368.800048828125; G=21373157; H=2024; I=11; J=1; K=0; L=0; M=0; N=44; O=0; P=0; Q=2},
    @{Row=69; A=45627; B=1361;                C=1563.349975585938;  D=1340;                E=1470.150024414062; G=26660938; H=2024; I=12; J=1; K=0; L=0; M=0; N=48; O=0; P=0; Q=0}
)

foreach ($r in $newRows) {
    $row = $r.Row

    # copy formatting (number formats, date style on col A, etc.) from the prior row
    $ws.Range("A" + ($row - 1) + ":R" + ($row - 1)).Copy()
    $ws.Range("A" + $row + ":R" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    # column F (Adj Close) and column R (backup) intentionally left blank for these new rows
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
}

Write-Host "done"
